$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.345.55'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '1.879.96'
$ws.Range('E3').Value = '  -1.59%  '
$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('D5').Value = "'246.65"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.73%  '
$ws.Range('D6').Value = "'0.689"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -5.46%  '
$ws.Range('E7').Value = '  -0.78%  '
$ws.Range('D8').Value = "'43.34"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +5.85%  '
$ws.Range('D9').Value = "'0.351"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.26%  '
$ws.Range('D10').Value = "'0.0738"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').Value = "'0.0969"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.04%  '
$ws.Range('D12').Value = "'13.13"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +1.32%  '
$ws.Range('D13').Value = '2.152.41'
$ws.Range('E13').Value = '  -1.57%  '
$ws.Range('D14').Value = "'0.737"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = "'4.95"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('D16').Value = '1.902.09'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '35.338.05'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = "'73.55"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.33%  '
$ws.Range('D19').Value = '0.0₃0822'
$ws.Range('E19').Value = '  -3.10%  '
$ws.Range('D20').Value = "'245.89"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = "'12.82"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.35%  '
$ws.Range('D22').Value = "'4.94"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.93%  '
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('D24').Value = "'2.57"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +4.90%  '
$ws.Range('E25').Value = '  -11.30%  '
$ws.Range('D26').Value = "'165.61"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('D27').Value = "'8.48"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.14%  '
$ws.Range('D28').Value = "'18.31"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('E29').Value = '  -4.33%  '
$ws.Range('D30').Value = '4.128.47'
$ws.Range('E31').Value = '  +4.93%  '
$ws.Range('D32').Value = "'4.24"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('D33').Value = "'0.0581"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.88%  '
$ws.Range('D34').Value = "'4.21"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.81%  '
$ws.Range('E35').Value = '  -0.76%  '
$ws.Range('D36').Value = "'0.852"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -6.98%  '
$ws.Range('E37').Value = '  -3.17%  '
$ws.Range('D38').Value = "'1.55"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -22.14%  '
$ws.Range('D39').Value = "'0.0691"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +7.07%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = "'97.34"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = "'17.01"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.80%  '
$ws.Range('E42').Value = '  -2.99%  '
$ws.Range('D43').Value = "'1.09"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.98%  '
$ws.Range('D44').Value = '1.290.23'
$ws.Range('E44').Value = '  -3.79%  '
$ws.Range('D45').Value = "'2.34"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -5.85%  '
$ws.Range('D46').Value = "'0.0810"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +7.13%  '
$ws.Range('E47').Value = '  -1.20%  '
$ws.Range('D48').Value = "'2.72"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').Value = "'12.13"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('D50').Value = "'43.26"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.21%  '
$ws.Range('D51').Value = "'6.26"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -7.13%  '

$wb.Save()
